$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.7521946666666667
$ws.Range("N2").Value = 2.256584
$ws.Range("O2").Value = 0.07361670343069449
$ws.Range("P2").Value = 0.0736167034306945
$ws.Range("Q2").Value = 0.01336198605866667
$ws.Range("R2").Value = 0.120257874528
$ws.Range("S2").Value = 0.07361670343069449
$ws.Range("T2").Value = 0.0736167034306945

# Row 3 updates
$ws.Range("O3").Value = 0.6908862423022597
$ws.Range("P3").Value = 0.6908862423022598
$ws.Range("S3").Value = 0.6908862423022597
$ws.Range("T3").Value = 0.6908862423022598

# Row 4 updates
$ws.Range("M4").Value = 2.406242333333334
$ws.Range("N4").Value = 7.218727
$ws.Range("O4").Value = 0.2354970542670457
$ws.Range("P4").Value = 0.2354970542670457
$ws.Range("S4").Value = 0.2354970542670457
$ws.Range("T4").Value = 0.2354970542670457
